# Refresh the crypto price/volume snapshot (columns D and E) for rows 2-51,
# matching the latest values pulled by the scheduled scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.614.02'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.882.45'
$ws.Range('D4').Value = "'0.9995"
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'246.56"
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = "'0.9996"
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').Value = "'0.4726"
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = "'0.2886"
$ws.Range('E8').Value = '  -1.19%  '
$ws.Range('D9').Value = "'0.06537"
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('D10').Value = "'22.19"
$ws.Range('E10').Value = '  +0.75%  '
$ws.Range('D11').Value = "'0.7762"
$ws.Range('E11').Value = '  +5.41%  '
$ws.Range('D12').Value = "'100.93"
$ws.Range('E12').Value = '  +4.49%  '
$ws.Range('D13').Value = "'0.07833"
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('D14').Value = '1.882.74'
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('D15').Value = "'5.253"
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = "'286.32"
$ws.Range('E16').Value = '  +1.02%  '
$ws.Range('D17').Value = '30.576.00'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = "'13.21"
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').Value = "'0.000007532"
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '2.126.54'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').Value = "'5.372"
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('D23').Value = "'0.9993"
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = "'6.411"
$ws.Range('E24').Value = '  +2.63%  '
$ws.Range('D25').Value = "'9.139"
$ws.Range('E25').Value = '  -0.98%  '
$ws.Range('D26').Value = "'162.90"
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('D27').Value = "'19.10"
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('D28').Value = "'1.917"
$ws.Range('D29').Value = "'0.09702"
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').Value = "'4.269"
$ws.Range('D33').Value = "'4.197"
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').Value = "'0.04855"
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('D35').Value = "'1.131"
$ws.Range('E35').Value = '  +0.34%  '
$ws.Range('D36').Value = "'0.6975"
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').Value = "'2.753"
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('D38').Value = "'0.01919"
$ws.Range('E38').Value = '  +1.37%  '
$ws.Range('D39').Value = "'2.881"
$ws.Range('E39').Value = '  +2.67%  '
$ws.Range('D40').Value = "'76.46"
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('D41').Value = "'6.294"
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('D42').Value = "'1.988"
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').Value = "'0.9993"
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').Value = "'0.8324"
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('D46').Value = "'101.50"
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').Value = "'9.810"
$ws.Range('E47').Value = '  +3.67%  '
$ws.Range('D48').Value = "'7.047"
$ws.Range('E48').Value = '  +0.29%  '
$ws.Range('D49').Value = "'35.21"
$ws.Range('E49').Value = '  -1.35%  '
$ws.Range('D50').Value = "'897.99"
$ws.Range('E50').Value = '  -2.03%  '
$ws.Range('D51').Value = "'0.05765"
$ws.Range('E51').Value = '  +0.20%  '
